$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 228.06667
$ws.Range("I41").Value = 290.7143
$ws.Range("J41").Value = 173.25
$ws.Range("K41").Value = 290.7143
$ws.Range("L41").Value = 173.25
$ws.Range("M41").Value = 149.2857
$ws.Range("N41").Value = -1053.25
$ws.Range("H62").Value = 6511
$ws.Range("I62").Value = 4866.6665
$ws.Range("K62").Value = 4866.6665
$ws.Range("M62").Value = -4242.6665
$ws.Range("H65").Value = 6511
$ws.Range("I65").Value = 4866.6665
$ws.Range("K65").Value = 24333.3325
$ws.Range("M65").Value = -21213.3325
$ws.Range("H86").Value = 3651.1333
$ws.Range("J86").Value = 4089.0833
$ws.Range("L86").Value = 4089.0833
$ws.Range("N86").Value = -6335.0833
$ws.Range("H89").Value = 3651.1333
$ws.Range("J89").Value = 4089.0833
$ws.Range("L89").Value = 20445.4165
$ws.Range("N89").Value = -31677.4165
$ws.Range("H112").Value = 2566.3333
$ws.Range("I112").Value = 774.5
$ws.Range("J112").Value = 2842
$ws.Range("K112").Value = 2323.5
$ws.Range("L112").Value = 8526
$ws.Range("M112").Value = -1215.5
$ws.Range("N112").Value = -10742
$ws.Range("H131").Value = 3962.913
$ws.Range("I131").Value = 1141.3572
$ws.Range("K131").Value = 3424.0716
$ws.Range("M131").Value = 1615.9284
$ws.Range("H135").Value = 831.8
$ws.Range("I135").Value = 621.56525
$ws.Range("J135").Value = 3249.5
$ws.Range("K135").Value = 5594.08725
$ws.Range("L135").Value = 29245.5
$ws.Range("M135").Value = -3059.08725
$ws.Range("N135").Value = -34315.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 27499.5
$ws.Range("I7").Value = 20000
$ws.Range("J7").Value = 34999
$ws.Range("K7").Value = 20000
$ws.Range("L7").Value = 34999
$ws.Range("M7").Value = -19886
$ws.Range("N7").Value = -35227
$ws.Range("H88").Value = 684
$ws.Range("I88").Value = 787.5
$ws.Range("J88").Value = 270
$ws.Range("K88").Value = 787.5
$ws.Range("L88").Value = 270
$ws.Range("M88").Value = -381.5
$ws.Range("N88").Value = -1082
$ws.Range("H91").Value = 684
$ws.Range("I91").Value = 787.5
$ws.Range("J91").Value = 270
$ws.Range("K91").Value = 787.5
$ws.Range("L91").Value = 270
$ws.Range("M91").Value = 616.5
$ws.Range("N91").Value = -3078
$ws.Range("H102").Value = 2750
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 2750
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 2750
$ws.Range("N102").Value = -5994
$ws.Range("M102").ClearContents()
$ws.Range("H122").Value = 590862.0600000001
$ws.Range("I122").Value = 1112720.6
$ws.Range("J122").Value = 3771.125
$ws.Range("K122").Value = 3338161.8
$ws.Range("L122").Value = 11313.375
$ws.Range("M122").Value = -3335711.8
$ws.Range("N122").Value = -16213.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1681.1177
$ws.Range("I20").Value = 815.6667
$ws.Range("J20").Value = 3758.2
$ws.Range("K20").Value = 815.6667
$ws.Range("L20").Value = 3758.2
$ws.Range("M20").Value = -568.6667
$ws.Range("N20").Value = -4252.2
$ws.Range("H86").Value = 3444
$ws.Range("I86").Value = 3612.2856
$ws.Range("J86").Value = 3275.7144
$ws.Range("K86").Value = 3612.2856
$ws.Range("L86").Value = 3275.7144
$ws.Range("M86").Value = -2489.2856
$ws.Range("N86").Value = -5521.7144
$ws.Range("H89").Value = 3444
$ws.Range("I89").Value = 3612.2856
$ws.Range("J89").Value = 3275.7144
$ws.Range("K89").Value = 18061.428
$ws.Range("L89").Value = 16378.572
$ws.Range("M89").Value = -12445.428
$ws.Range("N89").Value = -27610.572
$ws.Range("H94").Value = 6000
$ws.Range("I94").Value = 6000
$ws.Range("K94").Value = 6000
$ws.Range("M94").Value = -5549
$ws.Range("H99").Value = 4996.364
$ws.Range("I99").Value = 5096
$ws.Range("K99").Value = 5096
$ws.Range("M99").Value = -3598
$ws.Range("H105").Value = 4564.381
$ws.Range("I105").Value = 3520.7646
$ws.Range("J105").Value = 8999.75
$ws.Range("K105").Value = 3520.7646
$ws.Range("L105").Value = 8999.75
$ws.Range("M105").Value = -1773.7646
$ws.Range("N105").Value = -12493.75
$ws.Range("H134").Value = 2421.6924
$ws.Range("I134").Value = 2128.913
$ws.Range("K134").Value = 6386.739
$ws.Range("M134").Value = -3851.739

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1050.1111
$ws.Range("I16").Value = 956.375
$ws.Range("K16").Value = 956.375
$ws.Range("M16").Value = -669.375
$ws.Range("H31").Value = 3184.375
$ws.Range("J31").Value = 6579.1665
$ws.Range("L31").Value = 6579.1665
$ws.Range("N31").Value = -7169.1665
$ws.Range("H34").Value = 3184.375
$ws.Range("J34").Value = 6579.1665
$ws.Range("L34").Value = 6579.1665
$ws.Range("N34").Value = -6983.1665
$ws.Range("H99").Value = 11041.451
$ws.Range("I99").Value = 8593.546
$ws.Range("J99").Value = 12387.8
$ws.Range("K99").Value = 8593.546
$ws.Range("L99").Value = 12387.8
$ws.Range("M99").Value = -7095.546
$ws.Range("N99").Value = -15383.8
$ws.Range("H113").Value = 1050.1111
$ws.Range("I113").Value = 956.375
$ws.Range("K113").Value = 956.375
$ws.Range("M113").Value = 1213.625
$ws.Range("H122").Value = 4266.5806
$ws.Range("I122").Value = 4605.6313
$ws.Range("J122").Value = 3729.75
$ws.Range("K122").Value = 13816.8939
$ws.Range("L122").Value = 11189.25
$ws.Range("M122").Value = -11366.8939
$ws.Range("N122").Value = -16089.25
$ws.Range("H126").Value = 11041.451
$ws.Range("I126").Value = 8593.546
$ws.Range("J126").Value = 12387.8
$ws.Range("K126").Value = 25780.638
$ws.Range("L126").Value = 37163.39999999999
$ws.Range("M126").Value = -23310.638
$ws.Range("N126").Value = -42103.39999999999
$ws.Range("H134").Value = 3151.6428
$ws.Range("I134").Value = 2807
$ws.Range("J134").Value = 3611.1667
$ws.Range("K134").Value = 8421
$ws.Range("L134").Value = 10833.5001
$ws.Range("M134").Value = -5886
$ws.Range("N134").Value = -15903.5001
$ws.Range("H141").Value = 139165.33
$ws.Range("J141").Value = 139165.33
$ws.Range("L141").Value = 139165.33
$ws.Range("N141").Value = -149525.33

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 1163.3334
$ws.Range("J26").Value = 3000
$ws.Range("L26").Value = 9000
$ws.Range("N26").Value = -9576
$ws.Range("H32").Value = 6346428.5
$ws.Range("I32").Value = 1291.6666
$ws.Range("J32").Value = 13960593
$ws.Range("K32").Value = 3874.9998
$ws.Range("L32").Value = 41881779
$ws.Range("M32").Value = -3591.9998
$ws.Range("N32").Value = -41882345
$ws.Range("H39").Value = 3000
$ws.Range("I39").Value = 2400
$ws.Range("J39").Value = 4500
$ws.Range("K39").Value = 7200
$ws.Range("L39").Value = 13500
$ws.Range("M39").Value = -6906
$ws.Range("N39").Value = -14088
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H140").Value = 3280.625
$ws.Range("I140").Value = 3280.625
$ws.Range("K140").Value = 9841.875
$ws.Range("M140").Value = -4661.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 615.53845
$ws.Range("I2").Value = 183.57143
$ws.Range("K2").Value = 183.57143
$ws.Range("M2").Value = -70.57142999999999
$ws.Range("H80").Value = 1299.591
$ws.Range("I80").Value = 710.2857
$ws.Range("K80").Value = 710.2857
$ws.Range("M80").Value = 287.7143
$ws.Range("H83").Value = 1299.591
$ws.Range("I83").Value = 710.2857
$ws.Range("K83").Value = 3551.4285
$ws.Range("M83").Value = 1440.5715
$ws.Range("H132").Value = 2920.5833
$ws.Range("I132").Value = 2155.7144
$ws.Range("J132").Value = 3991.4
$ws.Range("K132").Value = 6467.1432
$ws.Range("L132").Value = 11974.2
$ws.Range("M132").Value = -3937.1432
$ws.Range("N132").Value = -17034.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3577.634
$ws.Range("I132").Value = 3465.0322
$ws.Range("J132").Value = 3926.7
$ws.Range("K132").Value = 10395.0966
$ws.Range("L132").Value = 11780.1
$ws.Range("M132").Value = -7865.096600000001
$ws.Range("N132").Value = -16840.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1130.9333
$ws.Range("I132").Value = 1144.1154
$ws.Range("K132").Value = 3432.3462
$ws.Range("M132").Value = -902.3462
$ws.Range("H133").Value = 80000
$ws.Range("J133").Value = 80000
$ws.Range("L133").Value = 80000
$ws.Range("N133").Value = -90120
$ws.Range("H136").Value = 2732.6667
$ws.Range("I136").Value = 1258.75
$ws.Range("J136").Value = 8628.333000000001
$ws.Range("K136").Value = 3776.25
$ws.Range("L136").Value = 25884.999
$ws.Range("M136").Value = -1226.25
$ws.Range("N136").Value = -30984.999
